# "Prova videos i captar el dia d'avui"
# - Shift every date in column B (rows 2-367) forward by 1050 days.
# - Add a new video-link cell F4 (reusing the existing shared string
#   "./video/tortura_1.mp4"), formatted like the other video-link cells
#   (same look as D2, underlined "link" style).
# - Update the sheet's active selection to C12 (and scroll so column B is
#   visible, matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every date in column B down by 1050 days (keeps existing s="4"
# date-number-format style; use Value2 so we operate on the raw serial
# number rather than the formatted display text).
for ($r = 2; $r -le 367; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 + 1050
}

# Add the new F4 video link cell, copying the underline "link" style used
# elsewhere in the sheet (e.g. D2) and then setting its text.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Value = "./video/tortura_1.mp4"

# Update view/selection to match what was saved: scroll to top-left B1 and
# select C12.
$excel.Goto($ws.Range("B1"), $true) | Out-Null
$ws.Range("C12").Select() | Out-Null
